# Apply the edit described by the diff:
# - Insert a new row (new row 2) above the existing "Clean Code..." row,
#   shifting it down to row 3.
# - New row 2 holds a new book entry with an empty/placeholder author & publisher list.
# - The old (now row 3) publisher value is corrected from
#   "['PReNtICE HALL']" to "['Prentice Hall Ptr']".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing the existing row 2 ("Clean Code...") to row 3.
$ws.Rows.Item(2).Insert()

# New row 2: newly added book entry.
$ws.Range("A2").Value = "MY BOOK  COVER Secrets in a Silicon Valley Startup"
$ws.Range("C2").Value = "[]"
$ws.Range("D2").Value = "[]"

# Row 3 (previously row 2) already retains Title/Author values after the insert;
# just fix the publisher value.
$ws.Range("D3").Value = "['Prentice Hall Ptr']"

$wb.Save()
